# Updates the "cryptos" price list on Sheet1 to the latest scraped
# snapshot. Two coin pairs swap rank order (rows 14/15 and 35/36) and
# two coins are replaced outright by newer ones that overtook them in
# the ranking (rows 39/40 and 42/43); every row's Price (D) and
# Volume(1h) (E) columns are refreshed with the new reading.
#
# Price values are text (e.g. "65.675.97", "0.998") rather than native
# numbers, so each one is written with a leading apostrophe to force
# Excel to keep it as text instead of re-parsing it as a number/date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'65.675.97"
$ws.Range("E2").Value = "  -5.83%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'3.266.09"
$ws.Range("E3").Value = "  -6.61%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'555.64"
$ws.Range("E5").Value = "  -4.19%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'184.50"
$ws.Range("E6").Value = "  -4.25%  "

# Row 7 - USDC
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.25%  "

# Row 8 - XRP
$ws.Range("D8").Value = "'0.587"
$ws.Range("E8").Value = "  -4.19%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "'3.257.67"
$ws.Range("E9").Value = "  -6.57%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.184"
$ws.Range("E10").Value = "  -10.24%  "

# Row 11 - Cardano
$ws.Range("D11").Value = "'0.581"
$ws.Range("E11").Value = "  -6.29%  "

# Row 12 - Avalanche
$ws.Range("D12").Value = "'46.90"
$ws.Range("E12").Value = "  -8.81%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -7.64%  "

# Row 14 - was Polkadot, now BitcoinCash (rows 14/15 swap order)
$ws.Range("B14").Value = "BitcoinCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D14").Value = "'634.10"
$ws.Range("E14").Value = "  -2.10%  "

# Row 15 - was BitcoinCash, now Polkadot
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'8.59"
$ws.Range("E15").Value = "  -6.08%  "

# Row 16 - WrappedliquidstakedEther2.0
$ws.Range("D16").Value = "'3.789.72"
$ws.Range("E16").Value = "  -6.36%  "

# Row 17 - Chainlink
$ws.Range("D17").Value = "'17.98"
$ws.Range("E17").Value = "  -1.77%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "'65.630.79"
$ws.Range("E18").Value = "  -5.84%  "

# Row 19 - TRON
$ws.Range("E19").Value = "  -3.37%  "

# Row 20 - WrappedEther
$ws.Range("D20").Value = "'3.257.20"
$ws.Range("E20").Value = "  -6.74%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -8.85%  "

# Row 22 - Polygon
$ws.Range("D22").Value = "'0.900"
$ws.Range("E22").Value = "  -5.17%  "

# Row 23 - InternetComputer(DFINITY)
$ws.Range("D23").Value = "'18.17"
$ws.Range("E23").Value = "  +0.52%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'106.67"
$ws.Range("E24").Value = "  +7.94%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "'4.89"
$ws.Range("E25").Value = "  -8.10%  "

# Row 26 - PancakeSwap
$ws.Range("D26").Value = "'3.95"
$ws.Range("E26").Value = "  -7.87%  "

# Row 27 - ImmutableX
$ws.Range("E27").Value = "  -7.68%  "

# Row 28 - RenderToken
$ws.Range("D28").Value = "'9.52"
$ws.Range("E28").Value = "  -5.19%  "

# Row 29 - Filecoin
$ws.Range("D29").Value = "'8.63"
$ws.Range("E29").Value = "  -7.86%  "

# Row 30 - EthereumClassic
$ws.Range("D30").Value = "'30.11"
$ws.Range("E30").Value = "  -7.90%  "

# Row 31 - dogwifhat
$ws.Range("D31").Value = "'3.94"
$ws.Range("E31").Value = "  -7.81%  "

# Row 32 - NEARProtocol
$ws.Range("D32").Value = "'6.23"
$ws.Range("E32").Value = "  -7.39%  "

# Row 33 - Cosmos
$ws.Range("D33").Value = "'11.00"
$ws.Range("E33").Value = "  -5.60%  "

# Row 34 - Hedera
$ws.Range("E34").Value = "  -5.03%  "

# Row 35 - was OKB, now Maker (rows 35/36 swap order)
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "'3.736.98"
$ws.Range("E35").Value = "  +0.74%  "

# Row 36 - was Maker, now OKB
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "'57.57"
$ws.Range("E36").Value = "  -5.93%  "

# Row 37 - Dai
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.22%  "

# Row 38 - Bittensor
$ws.Range("D38").Value = "'519.92"
$ws.Range("E38").Value = "  -6.16%  "

# Row 39 - was Stacks, now PEPE (rows 39/40 swap order)
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "'0.0₃0735"
$ws.Range("E39").Value = "  -7.09%  "

# Row 40 - was PEPE, now Stacks
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'3.36"
$ws.Range("E40").Value = "  -6.31%  "

# Row 41 - Kaspa
$ws.Range("E41").Value = "  -2.33%  "

# Row 42 - was Fetch.AI, now CoreDAO (rows 42/43 swap order)
$ws.Range("B42").Value = "CoreDAO"
$ws.Range("C42").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D42").Value = "'3.44"
$ws.Range("E42").Value = "  -5.77%  "

# Row 43 - was CoreDAO, now Fetch.AI
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").Value = "'2.70"
$ws.Range("E43").Value = "  -7.29%  "

# Row 44 - InjectiveProtocol
$ws.Range("D44").Value = "'32.81"
$ws.Range("E44").Value = "  -4.33%  "

# Row 45 - TheGraph
$ws.Range("D45").Value = "'0.335"
$ws.Range("E45").Value = "  -10.37%  "

# Row 46 - VeChain
$ws.Range("E46").Value = "  -7.31%  "

# Row 47 - ApeXProtocol
$ws.Range("D47").Value = "'3.20"
$ws.Range("E47").Value = "  -5.01%  "

# Row 48 - Stellar
$ws.Range("E48").Value = "  -4.33%  "

# Row 49 - ThetaToken
$ws.Range("D49").Value = "'2.59"
$ws.Range("E49").Value = "  -9.07%  "

# Row 50 - FirstDigitalUSD
$ws.Range("E50").Value = "  +0.09%  "

# Row 51 - Mantle
$ws.Range("D51").Value = "'1.25"
$ws.Range("E51").Value = "  +1.04%  "
